$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")
$ws.Range("C6").Value = 26.04
